$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "立讯精密"
$ws.Range("B3").Value = "东华软件"
$ws.Range("A4").Value = "福龙马"
$ws.Range("B4").Value = "立讯精密"
$ws.Range("C4").Value = "山子高科"
$ws.Range("A5").Value = "山子高科"
$ws.Range("B5").Value = "福龙马"
$ws.Range("C5").Value = "绝味食品"
$ws.Range("A6").Value = "绝味食品"
$ws.Range("B6").Value = "绝味食品"
$ws.Range("C6").Value = "凯美特气"
$ws.Range("A7").Value = "凯美特气"
$ws.Range("B7").Value = "天富能源"
$ws.Range("B8").Value = "天通股份"
$ws.Range("C8").Value = "福龙马"
$ws.Range("A9").Value = "东华软件"
$ws.Range("B9").Value = "联美控股"
$ws.Range("C9").Value = "上海建工"
$ws.Range("A10").Value = "山河智能"
$ws.Range("B10").Value = "赣锋锂业"
$ws.Range("C10").Value = "欧菲光"
$ws.Range("A11").Value = "上海建工"
$ws.Range("B11").Value = "凯美特气"
$ws.Range("C11").Value = "山河智能"
$ws.Range("A12").Value = "杭电股份"
$ws.Range("B12").Value = "山子高科"
$ws.Range("C12").Value = "首开股份"
$ws.Range("A13").Value = "天通股份"
$ws.Range("B13").Value = "山河智能"
$ws.Range("C13").Value = "赣锋锂业"
$ws.Range("A14").Value = "卧龙电驱"
$ws.Range("B14").Value = "杭电股份"
$ws.Range("A15").Value = "欧菲光"
$ws.Range("B15").Value = "东方财富"
$ws.Range("C15").Value = "天普股份"
$ws.Range("A16").Value = "天富能源"
$ws.Range("B16").Value = "上海建工"
$ws.Range("C16").Value = "杭电股份"
$ws.Range("A17").Value = "联美控股"
$ws.Range("B17").Value = "长飞光纤"
$ws.Range("A18").Value = "省广集团"
$ws.Range("B18").Value = "中科通达"
$ws.Range("C18").Value = "均胜电子"
$ws.Range("A19").Value = "润和软件"
$ws.Range("B19").Value = "欧菲光"
$ws.Range("C19").Value = "金发科技"
$ws.Range("A20").Value = "歌尔股份"
$ws.Range("B20").Value = "万马股份"
$ws.Range("C20").Value = "三花智控"
$ws.Range("A21").Value = "首开股份"
$ws.Range("B21").Value = "卧龙电驱"
$ws.Range("C21").Value = "华胜天成"
